$d = $word.ActiveDocument
$bullet = [char]0x2022

# -----------------------------------------------------------------------
# This paragraph (PINCHme / Data Analyst bullets) gets two edits:
#
#  1) "<bullet> Pulled relevant and impactful data, metrics, and trends"
#     is split into five runs:
#        "<bullet> "
#        "Developed ad hoc SQL reports"
#        " in Treasure Data"
#        " to pull"
#        " relevant data, metrics, and trends"
#
#  2) A new run "and email list " is inserted right after the existing
#     "data pull " run (before "requests that came in...").
#
# The underlying engine coalesces same-formatted runs whenever a
# paragraph is touched, so after performing both text edits we walk the
# whole paragraph and re-assert every intended run boundary by toggling
# a character property (Bold on, then back off) across each exact
# sub-range. That forces a run break at both ends of the sub-range
# without altering the visible formatting.
# -----------------------------------------------------------------------

function Force-RunBoundaries($ranges) {
    foreach ($span in $ranges) {
        $r = $d.Range($span[0], $span[1])
        $r.Font.Bold = $true
        $r.Font.Bold = $false
    }
}

# --- locate the target paragraph text -----------------------------------
$find1 = $d.Content
$find1.Find.Execute(($bullet + " Pulled relevant and impactful data, metrics, and trends"), `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$run0Start = $find1.Start

# --- Change 1: rewrite the first run's text into 5 runs ------------------
$newParts = @(
    ($bullet + " "),
    "Developed ad hoc SQL reports",
    " in Treasure Data",
    " to pull",
    " relevant data, metrics, and trends"
)
$find1.Text = [string]::Join("", $newParts)

# --- Change 2: find the "data pull " run and insert a new run after it --
$find2 = $d.Content
$find2.Find.Execute("Completed data pull requests that came in from all branches within the business", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$dataPullEnd = $find2.Start + "Completed ".Length + "data pull ".Length
$insertion = "and email list "
$insRng = $d.Range($dataPullEnd, $dataPullEnd)
$insRng.InsertBefore($insertion)

# -----------------------------------------------------------------------
# Re-derive every run boundary in the whole paragraph (post edits) and
# force it, so previously-distinct runs don't get merged back together
# and the newly created runs stay split apart.
# -----------------------------------------------------------------------
$allParts = @()
$allParts += $newParts
$allParts += " from multimillion row database"
$allParts += "s"
$allParts += " to analyze and translate "
$allParts += "this"
$allParts += " into actionable insights for the business."
$allParts += "`n"
$allParts += "`n$bullet "
$allParts += "Facilitated"
$allParts += " a data-driven culture by building"
$allParts += " and managing several"
$allParts += " dashboards to display "
$allParts += "impactful"
$allParts += " metrics to our CEO, "
$allParts += "CFO,"
$allParts += " and other key stakeholders"
$allParts += " within the company."
$allParts += "`n"
$allParts += "`n$bullet "
$allParts += "Completed "
$allParts += "data pull "
$allParts += $insertion
$allParts += "requests that came in from all branches within the business"
$allParts += "."

$spans = New-Object System.Collections.ArrayList
$p = $run0Start
foreach ($part in $allParts) {
    [void]$spans.Add(@($p, $p + $part.Length))
    $p += $part.Length
}

Force-RunBoundaries $spans

"edit complete"
